$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_G"
$ws.Range("A2").Value = 73.354661791590487
$ws.Range("A3").Value = 72.577696526508234
$ws.Range("A4").Value = 71.023765996343684
$ws.Range("A5").Value = 71.435100548446073
$ws.Range("A6").Value = 71.572212065813517
$ws.Range("A7").Value = 72.760511882998173
$ws.Range("A8").Value = 70.42961608775137
$ws.Range("A9").Value = 71.297989031078615
$ws.Range("A10").Value = 70.155393053016454
$ws.Range("A11").Value = 70.292504570383912
$ws.Range("A12").Value = 74.908592321755023
$ws.Range("A13").Value = 75.365630712979893
$ws.Range("A14").Value = 70.338208409506393
$ws.Range("A15").Value = 71.06946983546618
$ws.Range("A16").Value = 70.749542961608782
$ws.Range("A17").Value = 73.354661791590487
$ws.Range("A18").Value = 74.954296160877504
$ws.Range("A19").Value = 75.045703839122496
$ws.Range("A20").Value = 73.171846435100548
$ws.Range("A21").Value = 70.886654478976226
$ws.Range("A22").Value = 73.720292504570381
$ws.Range("A23").Value = 71.755027422303471
$ws.Range("A24").Value = 73.628884826325418
$ws.Range("A25").Value = 73.400365630712983
$ws.Range("A26").Value = 71.572212065813517
$ws.Range("A27").Value = 71.480804387568554
$ws.Range("A28").Value = 71.206581352833638
$ws.Range("A29").Value = 72.669104204753197
$ws.Range("A30").Value = 71.755027422303471
$ws.Range("A31").Value = 71.846435100548447
$ws.Range("A32").Value = 70.338208409506393
$ws.Range("A33").Value = 70.383912248628889
$ws.Range("A34").Value = 70.475319926873851
$ws.Range("A35").Value = 71.846435100548447
$ws.Range("A36").Value = 72.166361974405845
$ws.Range("A37").Value = 76.279707495429619
$ws.Range("A38").Value = 71.206581352833638
$ws.Range("A39").Value = 71.572212065813517
$ws.Range("A40").Value = 72.166361974405845
$ws.Range("A41").Value = 71.709323583180989
$ws.Range("A42").Value = 72.212065813528341
$ws.Range("A43").Value = 71.93784277879341
$ws.Range("A44").Value = 71.663619744058508
$ws.Range("A45").Value = 72.257769652650822
$ws.Range("A46").Value = 70.749542961608782
$ws.Range("A47").Value = 70.749542961608782
$ws.Range("A48").Value = 72.943327239488127
$ws.Range("A49").Value = 71.06946983546618
